$wb = $excel.ActiveWorkbook

# --- Sheet "answers" (sheet1) ---
$ws1 = $wb.Worksheets.Item("answers")

# Header row
$ws1.Range("A1").Value = "Question"
$ws1.Range("B1").Value = "Student"
$ws1.Range("C1").Value = "Score"

# Data rows
$ws1.Range("A2").Value = "Question 1"
$ws1.Range("B2").Value = "Student 1"
$ws1.Range("C2").Value = 1

$ws1.Range("A3").Value = "Question 1"
$ws1.Range("B3").Value = "Student 2"
$ws1.Range("C3").Value = 1

$ws1.Range("A4").Value = "Question 2"
$ws1.Range("B4").Value = "Student 1"
$ws1.Range("C4").Value = 0

$ws1.Range("A5").Value = "Question 2"
$ws1.Range("B5").Value = "Student 2"
$ws1.Range("C5").Value = 0

$ws1.Range("A6").Value = "Question 3"
$ws1.Range("B6").Value = "Student 1"
$ws1.Range("C6").Value = 1

$ws1.Range("A7").Value = "Question 3"
$ws1.Range("B7").Value = "Student 2"
$ws1.Range("C7").Value = 0

$ws1.Range("B8").Select()

# --- Sheet "metadata" (sheet2) ---
$ws2 = $wb.Worksheets.Item("metadata")

$ws2.Range("B2").Value = "answers table contains the results of of an exam. The first row contains the question, the second contains the student name, and the third one contains the result for each question. 1 means they provided the right answer, 0 means they provided the wrong answer"

$ws2.Range("B2").Select()
